$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update UnitPrice (column G) values with a small surcharge applied to each row.
$ws.Range("G2").Value = 8.75
$ws.Range("G3").Value = 15.35
$ws.Range("G4").Value = 5.25
$ws.Range("G5").Value = 75.99
$ws.Range("G6").Value = 8.75
$ws.Range("G7").Value = 15.35
$ws.Range("G8").Value = 5.25
$ws.Range("G9").Value = 8.75
$ws.Range("G10").Value = 15.35
$ws.Range("G11").Value = 5.25
$ws.Range("G12").Value = 8.75
$ws.Range("G13").Value = 15.35
$ws.Range("G14").Value = 5.25
$ws.Range("G15").Value = 8.75
$ws.Range("G16").Value = 8.75
$ws.Range("G17").Value = 15.35
$ws.Range("G18").Value = 5.25
$ws.Range("G19").Value = 60
$ws.Range("G20").Value = 5.25
$ws.Range("G21").Value = 75.99
$ws.Range("G22").Value = 5.25
$ws.Range("G23").Value = 75.99
$ws.Range("G24").Value = 5.25
$ws.Range("G25").Value = 75.99

# NumUnits (column H) corrections
$ws.Range("H11").Value = 7
$ws.Range("H13").Value = 2

# Row 19 ProductID corrected from LMIPA to TPCA
$ws.Range("B19").Value = "TPCA"

# Autofit the now-used UnitPrice column to its new (longer, decimal) values
$ws.Columns("G:G").AutoFit()

# Restore the view state (scroll position and selection) as last left by the user
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("K22").Select()
